$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell E8 used to read "Good Morning"; the author updated it to "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Leave the active cell/selection on E8, matching the saved view state.
$ws.Range("E8").Select()
